$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38; this shifts existing rows 38-90 down to 39-91
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new data record
$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(45125)
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 100112042
$ws.Range("G38").Value = "Locoto"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 60
$ws.Range("K38").Value = 3800
$ws.Range("L38").Value = 3800
$ws.Range("M38").Value = 3800
$ws.Range("N38").Value = '$/kilo'
$ws.Range("O38").Value = "Región de Arica y Parinacota"
$ws.Range("P38").Value = 3800
$ws.Range("Q38").Value = 1
$ws.Range("R38").Value = "Hortaliza"
